# Apply the 2026-02-18 Betfair odds update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Value tweaks on the existing rows 2-6 (Turkish 1 Lig / UEFA CL / La Liga)
# ---------------------------------------------------------------------------

# Row 2 - Keciorengucu vs Erzurum BB
$ws.Cells.Item(2, 9).Value  = 2.92    # I2
$ws.Cells.Item(2, 12).Value = 1.32    # L2
$ws.Cells.Item(2, 20).Value = 1.73    # T2
$ws.Cells.Item(2, 21).Value = 2.08    # U2
$ws.Cells.Item(2, 22).Value = 1.52    # V2
$ws.Cells.Item(2, 34).Value = 970     # AH2

# Row 3 - Van Buyuksehir Belediye vs Bodrum Belediyesi Bodru
$ws.Cells.Item(3, 12).Value = 1.01    # L3
$ws.Cells.Item(3, 13).Value = 1.01    # M3
$ws.Cells.Item(3, 14).Value = 1.8     # N3
$ws.Cells.Item(3, 15).Value = 1.01    # O3
$ws.Cells.Item(3, 16).Value = 1.8     # P3
$ws.Cells.Item(3, 17).Value = 1.73    # Q3
$ws.Cells.Item(3, 18).Value = 1.05    # R3
$ws.Cells.Item(3, 19).Value = 1.73    # S3
$ws.Cells.Item(3, 20).Value = 1.01    # T3
$ws.Cells.Item(3, 21).Value = 1.01    # U3
$ws.Cells.Item(3, 22).Value = 1.46    # V3
$ws.Cells.Item(3, 23).Value = 1.44    # W3
$ws.Cells.Item(3, 24).Value = 1000    # X3
$ws.Cells.Item(3, 25).Value = 1000    # Y3
$ws.Cells.Item(3, 26).Value = 1000    # Z3
$ws.Cells.Item(3, 27).Value = 1000    # AA3
$ws.Cells.Item(3, 28).Value = 1000    # AB3
$ws.Cells.Item(3, 29).Value = 1000    # AC3
$ws.Cells.Item(3, 30).Value = 1000    # AD3
$ws.Cells.Item(3, 31).Value = 1000    # AE3
$ws.Cells.Item(3, 32).Value = 1000    # AF3
$ws.Cells.Item(3, 33).Value = 1000    # AG3
$ws.Cells.Item(3, 34).Value = 1000    # AH3
$ws.Cells.Item(3, 35).Value = 1000    # AI3
$ws.Cells.Item(3, 36).Value = 1000    # AJ3
$ws.Cells.Item(3, 37).Value = 1000    # AK3
$ws.Cells.Item(3, 38).Value = 1000    # AL3
$ws.Cells.Item(3, 39).Value = 1000    # AM3
$ws.Cells.Item(3, 40).Value = 1000    # AN3
$ws.Cells.Item(3, 41).Value = 1000    # AO3

# Row 4 - Manisa FK vs Bandirmaspor
$ws.Cells.Item(4, 6).Value  = 2.34    # F4
$ws.Cells.Item(4, 7).Value  = 2.66    # G4
$ws.Cells.Item(4, 8).Value  = 2.8     # H4
$ws.Cells.Item(4, 9).Value  = 3.55    # I4
$ws.Cells.Item(4, 10).Value = 3.45    # J4
$ws.Cells.Item(4, 11).Value = 4.1     # K4
$ws.Cells.Item(4, 16).Value = 2.06    # P4
$ws.Cells.Item(4, 17).Value = 1.73    # Q4

# Row 5 - Qarabag FK vs Newcastle
$ws.Cells.Item(5, 7).Value  = 9       # G5
$ws.Cells.Item(5, 8).Value  = 1.43    # H5
$ws.Cells.Item(5, 9).Value  = 1.44    # I5
$ws.Cells.Item(5, 10).Value = 5.2     # J5
$ws.Cells.Item(5, 15).Value = 1.23    # O5
$ws.Cells.Item(5, 16).Value = 2.36    # P5
$ws.Cells.Item(5, 17).Value = 1.69    # Q5
$ws.Cells.Item(5, 19).Value = 2.72    # S5
$ws.Cells.Item(5, 20).Value = 1.96    # T5
$ws.Cells.Item(5, 21).Value = 2.02    # U5
$ws.Cells.Item(5, 27).Value = 12      # AA5
$ws.Cells.Item(5, 28).Value = 30      # AB5
$ws.Cells.Item(5, 29).Value = 12      # AC5
$ws.Cells.Item(5, 31).Value = 14.5    # AE5
$ws.Cells.Item(5, 34).Value = 25      # AH5
$ws.Cells.Item(5, 36).Value = 320     # AJ5
$ws.Cells.Item(5, 37).Value = 140     # AK5
$ws.Cells.Item(5, 38).Value = 110     # AL5
$ws.Cells.Item(5, 39).Value = 150     # AM5
$ws.Cells.Item(5, 40).Value = 170     # AN5
$ws.Cells.Item(5, 41).Value = 5.9     # AO5

# Row 6 - Levante vs Villarreal
$ws.Cells.Item(6, 8).Value  = 1.89    # H6
$ws.Cells.Item(6, 13).Value = 1.05    # M6
$ws.Cells.Item(6, 14).Value = 4.8     # N6
$ws.Cells.Item(6, 17).Value = 1.73    # Q6
$ws.Cells.Item(6, 19).Value = 2.82    # S6
$ws.Cells.Item(6, 20).Value = 1.7     # T6
$ws.Cells.Item(6, 21).Value = 2.36    # U6
$ws.Cells.Item(6, 24).Value = 20      # X6
$ws.Cells.Item(6, 28).Value = 19      # AB6
$ws.Cells.Item(6, 29).Value = 9.4     # AC6
$ws.Cells.Item(6, 32).Value = 34      # AF6
$ws.Cells.Item(6, 33).Value = 17.5    # AG6
$ws.Cells.Item(6, 36).Value = 100     # AJ6
$ws.Cells.Item(6, 37).Value = 980     # AK6
$ws.Cells.Item(6, 39).Value = 80      # AM6

# ---------------------------------------------------------------------------
# 2) Insert a brand-new fixture row at row 8 (Italian Serie A, AC Milan v
#    Como), which pushes the former rows 8-11 down to rows 9-12.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Insert()

$newRow8 = @{
    1  = "Italian Serie A"
    2  = "'2026-02-18"
    3  = "16:45:00"
    4  = "AC Milan"
    5  = "Como"
    6  = 2.3
    7  = 2.38
    8  = 3.4
    9  = 3.5
    10 = 3.4
    11 = 3.55
    12 = 0
    13 = 1.07
    14 = 3.75
    15 = 1.33
    16 = 1.92
    17 = 2.02
    18 = 1.35
    19 = 3.6
    20 = 1.8
    21 = 2.16
    22 = 0
    23 = 0
    24 = 14
    25 = 13.5
    26 = 26
    27 = 65
    28 = 10.5
    29 = 7.8
    30 = 15
    31 = 42
    32 = 15
    33 = 11
    34 = 18
    35 = 55
    36 = 32
    37 = 27
    38 = 40
    39 = 110
    40 = 19.5
    41 = 42
}
foreach ($col in $newRow8.Keys) {
    $ws.Cells.Item(8, $col).Value = $newRow8[$col]
}

# ---------------------------------------------------------------------------
# 3) Value tweaks on the rows that shifted down (now rows 9-12)
# ---------------------------------------------------------------------------

# Row 9 (was row 8) - Wolves vs Arsenal
$ws.Cells.Item(9, 14).Value = 4.8     # N9
$ws.Cells.Item(9, 15).Value = 1.25    # O9
$ws.Cells.Item(9, 16).Value = 2.28    # P9
$ws.Cells.Item(9, 17).Value = 1.75    # Q9
$ws.Cells.Item(9, 18).Value = 1.51    # R9
$ws.Cells.Item(9, 26).Value = 7.2     # Z9
$ws.Cells.Item(9, 28).Value = 38      # AB9
$ws.Cells.Item(9, 29).Value = 14      # AC9
$ws.Cells.Item(9, 31).Value = 14      # AE9
$ws.Cells.Item(9, 33).Value = 55      # AG9
$ws.Cells.Item(9, 39).Value = 420     # AM9
$ws.Cells.Item(9, 41).Value = 5.1     # AO9

# Row 10 (was row 9) - Olympiakos vs Leverkusen
$ws.Cells.Item(10, 7).Value  = 2.7    # G10
$ws.Cells.Item(10, 8).Value  = 2.9    # H10
$ws.Cells.Item(10, 14).Value = 3.9    # N10
$ws.Cells.Item(10, 20).Value = 1.75   # T10
$ws.Cells.Item(10, 25).Value = 12     # Y10

# Row 11 (was row 10) - Club Brugge vs Atletico Madrid
$ws.Cells.Item(11, 6).Value  = 4      # F11
$ws.Cells.Item(11, 8).Value  = 1.98   # H11
$ws.Cells.Item(11, 9).Value  = 2      # I11
$ws.Cells.Item(11, 14).Value = 5      # N11
$ws.Cells.Item(11, 15).Value = 1.23   # O11
$ws.Cells.Item(11, 16).Value = 2.38   # P11
$ws.Cells.Item(11, 18).Value = 1.55   # R11
$ws.Cells.Item(11, 19).Value = 2.7    # S11
$ws.Cells.Item(11, 20).Value = 1.65   # T11
$ws.Cells.Item(11, 21).Value = 2.44   # U11
$ws.Cells.Item(11, 24).Value = 21     # X11
$ws.Cells.Item(11, 28).Value = 20     # AB11
$ws.Cells.Item(11, 29).Value = 9.4    # AC11
$ws.Cells.Item(11, 30).Value = 10.5   # AD11
$ws.Cells.Item(11, 31).Value = 19.5   # AE11
$ws.Cells.Item(11, 32).Value = 34     # AF11
$ws.Cells.Item(11, 33).Value = 17.5   # AG11
$ws.Cells.Item(11, 37).Value = 46     # AK11
$ws.Cells.Item(11, 38).Value = 48     # AL11
$ws.Cells.Item(11, 40).Value = 36     # AN11
$ws.Cells.Item(11, 41).Value = 10     # AO11

# Row 12 (was row 11) - Bodo Glimt vs Inter
$ws.Cells.Item(12, 7).Value  = 5.2    # G12
$ws.Cells.Item(12, 8).Value  = 1.71   # H12
$ws.Cells.Item(12, 9).Value  = 1.72   # I12
$ws.Cells.Item(12, 16).Value = 2.7    # P12
$ws.Cells.Item(12, 17).Value = 1.55   # Q12
$ws.Cells.Item(12, 18).Value = 1.69   # R12
$ws.Cells.Item(12, 19).Value = 2.36   # S12
$ws.Cells.Item(12, 20).Value = 1.61   # T12
$ws.Cells.Item(12, 25).Value = 13     # Y12
$ws.Cells.Item(12, 27).Value = 18.5   # AA12
$ws.Cells.Item(12, 32).Value = 50     # AF12
$ws.Cells.Item(12, 39).Value = 65     # AM12

Write-Output "edit applied"
